# Lab Inventory update:
#  - Insert a new row at row 54 (the blank spacer row right before the
#    "Tape" section header) and fill it with a new "Table Clamp" entry.
#    This pushes every row from the old row 54 onward down by one.
#  - Append a new "IR Sensor cables" entry as the final row of the sheet,
#    including a link to the part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Table Clamp" row ------------------------------------------
# Row 54 is currently an empty spacer row; inserting here shifts the
# existing rows 54-109 down to 55-110 and leaves a fresh blank row 54
# for the new item.
$ws.Rows(54).Insert()

$ws.Cells.Item(54, 1).Value2 = "Table Clamp"
$ws.Cells.Item(54, 2).Value2 = 1
$ws.Cells.Item(54, 3).Value2 = 1

# --- Append "IR Sensor cables" row -------------------------------------
$ws.Cells.Item(111, 1).Value2 = "IR Sensor cables"
$ws.Cells.Item(111, 2).Value2 = 12
$ws.Cells.Item(111, 3).Value2 = 100
$ws.Cells.Item(111, 4).Value2 = "https://www.pololu.com/product/117"

# Match the workbook's saved selection/view state from the edit.
$ws.Range("D111").Select()
